$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before column C ("descriptions" column) ---
# This shifts the old C..AH columns to D..AI, and shifts the existing
# merged cell ranges along with them automatically.
$ws.Columns("C").Insert()

# --- 2. Add the three new header rows above the register table ---
$ws.Range("A1").Value = "Module name"
$ws.Range("B1").Value = "gpio_regfile"

$ws.Range("A2").Value = "Base Addr"
$ws.Range("B2").Value = "0x40000000"

$ws.Range("A3").Value = "Data Width"
$ws.Range("B3").Value = 32

# --- 3. Populate the new "descriptions" column ---
$ws.Range("C4").Value = "descriptions"

# Blank, merged description cells next to each register block.
$ws.Range("C5:C7").Merge()
$ws.Range("C8:C10").Merge()
$ws.Range("C11:C13").Merge()
$ws.Range("C14:C16").Merge()

# --- 4. Fix up the label in the last row of the table ---
$ws.Range("B16").Value = " Attributes"

# --- 5. Column widths ---
$ws.Columns("A").ColumnWidth = 12.51
$ws.Columns("B").ColumnWidth = 12.36
$ws.Columns("C").ColumnWidth = 12.36

# --- 6. Selection / view niceties to mirror the saved workbook state ---
$ws.Range("B4").Select()
